$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.668.00"
$ws.Range("E2").Value = "  -3.40%  "
$ws.Range("D3").Value = "3.170.94"
$ws.Range("E3").Value = "  -7.90%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").Value = "'565.14"
$ws.Range("E5").Value = "  -3.70%  "
$ws.Range("D6").Value = "'170.73"
$ws.Range("E6").Value = "  -5.08%  "
$ws.Range("D7").Value = "'0.615"
$ws.Range("E7").Value = "  -2.52%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("D9").Value = "3.167.50"
$ws.Range("E9").Value = "  -7.99%  "
$ws.Range("D10").Value = "'0.125"
$ws.Range("E10").Value = "  -6.63%  "
$ws.Range("D11").Value = "'6.57"
$ws.Range("E11").Value = "  -5.68%  "
$ws.Range("D12").Value = "'0.395"
$ws.Range("E12").Value = "  -5.43%  "
$ws.Range("D13").Value = "3.719.74"
$ws.Range("E13").Value = "  -8.08%  "
$ws.Range("E14").Value = "  +0.57%  "
$ws.Range("D15").Value = "'27.38"
$ws.Range("E15").Value = "  -8.97%  "
$ws.Range("D16").Value = "64.551.36"
$ws.Range("E16").Value = "  -3.51%  "
$ws.Range("D17").Value = "'0.0000163"
$ws.Range("E17").Value = "  -6.62%  "
$ws.Range("D18").Value = "3.158.39"
$ws.Range("E18").Value = "  -8.28%  "
$ws.Range("D19").Value = "'5.74"
$ws.Range("E19").Value = "  -3.78%  "
$ws.Range("D20").Value = "'13.01"
$ws.Range("E20").Value = "  -6.41%  "
$ws.Range("D21").Value = "'354.96"
$ws.Range("E21").Value = "  -4.93%  "
$ws.Range("D22").Value = "'7.26"
$ws.Range("E22").Value = "  -5.51%  "
$ws.Range("D23").Value = "'1.01"
$ws.Range("E23").Value = "  +0.79%  "
$ws.Range("D24").Value = "'68.42"
$ws.Range("E24").Value = "  -7.03%  "
$ws.Range("D25").Value = "'0.502"
$ws.Range("E25").Value = "  -6.57%  "
$ws.Range("D26").Value = "'0.0000118"
$ws.Range("E26").Value = "  -9.91%  "
$ws.Range("D27").Value = "'9.66"
$ws.Range("E27").Value = "  -3.58%  "
$ws.Range("E28").Value = "  -2.09%  "
$ws.Range("D29").Value = "'0.998"
$ws.Range("E29").Value = "  -0.16%  "
$ws.Range("E30").Value = "  -0.20%  "
$ws.Range("D31").Value = "'1.90"
$ws.Range("E31").Value = "  -5.24%  "
$ws.Range("D32").Value = "'5.41"
$ws.Range("E32").Value = "  -8.06%  "
$ws.Range("D33").Value = "'22.01"
$ws.Range("E33").Value = "  -7.12%  "
$ws.Range("E34").Value = "  -6.03%  "
$ws.Range("D35").Value = "'6.67"
$ws.Range("E35").Value = "  -6.39%  "
$ws.Range("D36").Value = "'1.44"
$ws.Range("E36").Value = "  -8.87%  "
$ws.Range("D37").Value = "'153.57"
$ws.Range("E37").Value = "  -5.81%  "
$ws.Range("D38").Value = "'0.827"
$ws.Range("E38").Value = "  -6.27%  "
$ws.Range("D39").Value = "'26.16"
$ws.Range("E39").Value = "  -6.18%  "
$ws.Range("D40").Value = "'1.73"
$ws.Range("E40").Value = "  -5.16%  "
$ws.Range("D41").Value = "'2.51"
$ws.Range("E41").Value = "  -5.53%  "
$ws.Range("D42").Value = "2.651.61"
$ws.Range("E42").Value = "  -3.68%  "
$ws.Range("D43").Value = "'4.18"
$ws.Range("E43").Value = "  -7.19%  "
$ws.Range("B44").Value = "RenderToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D44").Value = "'6.04"
$ws.Range("E44").Value = "  -5.28%  "
$ws.Range("B45").Value = "OKB"
$ws.Range("C45").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D45").Value = "'39.42"
$ws.Range("E45").Value = "  -1.98%  "
$ws.Range("B46").Value = "InjectiveProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D46").Value = "'24.13"
$ws.Range("E46").Value = "  -6.01%  "
$ws.Range("B47").Value = "Hedera"
$ws.Range("C47").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D47").Value = "'0.0653"
$ws.Range("E47").Value = "  -6.74%  "
$ws.Range("D48").Value = "'319.24"
$ws.Range("E48").Value = "  -5.97%  "
$ws.Range("D49").Value = "'0.0273"
$ws.Range("E49").Value = "  -5.02%  "
$ws.Range("D50").Value = "'0.102"
$ws.Range("E50").Value = "  -3.51%  "
$ws.Range("D51").Value = "'0.999"
$ws.Range("E51").Value = "  -0.17%  "
